$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.003.09"
$ws.Range("E2").Value = "  +0.42%  "

$ws.Range("D3").Value = "1.561.56"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "207.48"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").Value = "0.490"
$ws.Range("E6").Value = "  +1.16%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "22.16"
$ws.Range("E8").Value = "  +2.11%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").Value = "0.0597"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").Value = "0.0860"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").Value = "1.784.75"
$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("D13").Value = "1.544.45"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").Value = "62.10"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").Value = "27.000.85"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("D19").Value = "217.10"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("E20").Value = "  +2.22%  "

$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("E22").Value = "  +1.56%  "

$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("E24").Value = "  -2.53%  "

$ws.Range("D25").Value = "153.17"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "6.63"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +1.16%  "

$ws.Range("E28").Value = "  +1.39%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  +2.03%  "

$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("D33").Value = "1.424.27"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("E34").Value = "  +3.81%  "

$ws.Range("E35").Value = "  +3.24%  "

$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  +9.11%  "

$ws.Range("E37").Value = "  +1.26%  "

$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("E39").Value = "  +1.86%  "

$ws.Range("D40").Value = "0.810"
$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").Value = "5.68"
$ws.Range("E42").Value = "  +0.75%  "

$ws.Range("E43").Value = "  +2.98%  "

$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").Value = "64.98"
$ws.Range("E45").Value = "  +2.11%  "

$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("D47").Value = "1.697.44"
$ws.Range("E47").Value = "  +0.65%  "

$ws.Range("D48").Value = "87.54"
$ws.Range("E48").Value = "  +1.58%  "

$ws.Range("D49").Value = "0.0522"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("E50").Value = "  +0.50%  "

$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  -0.40%  "
